# Rename the two logo pictures that live in the document's headers/footers.
#
#   * BTec_Logo-Orange  (the two headers)   : image2.jpg -> image1.jpg
#   * PearsonLogo       (the two footers)   : image1.png -> image2.png
#
# The picture "name" isn't changed anywhere in the body text - it is the
# name stamped on the drawing object itself (wp:docPr/@name, mirrored on
# pic:cNvPr/@name), so it has to be touched through each InlineShape's
# .Name property rather than Find/Replace (which only sees run text).

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
$headerPrimary = $sec.Headers.Item(1)
$headerFirst   = $sec.Headers.Item(2)
$footerPrimary = $sec.Footers.Item(1)
$footerFirst   = $sec.Footers.Item(2)

# --- Headers: BTec_Logo-Orange, image2.jpg -> image1.jpg ---
$hPrimaryShape = $headerPrimary.Range.InlineShapes.Item(1)
$hPrimaryShape.Range.InlineShapes.Item(1).Name = "image1.jpg"

$hFirstShape = $headerFirst.Range.InlineShapes.Item(1)
$hFirstShape.Range.InlineShapes.Item(1).Name = "image1.jpg"

# --- Footers: PearsonLogo, image1.png -> image2.png ---
$fPrimaryShape = $footerPrimary.Range.InlineShapes.Item(1)
$fPrimaryShape.Range.InlineShapes.Item(1).Name = "image2.png"

$fFirstShape = $footerFirst.Range.InlineShapes.Item(1)
$fFirstShape.Range.InlineShapes.Item(1).Name = "image2.png"

Write-Host "Renamed header/footer logo pictures"
